$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.915.69'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.863.37'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5048'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3648'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07186'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8946'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').Value = '1.888.91'
$ws.Range('E12').Value = '  +1.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07499'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.237'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.0000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008513'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '26.961.55'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.028'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').Value = '2.122.32'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.413'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.783'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.093'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.707'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.676'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09224'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05144'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7507'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.957'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.155'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.248'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.600'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02003'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5581'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.070'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.575'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.47'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.579'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.43%  '
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4704'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.12%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.05'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.563'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.12'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.25%  '
